# Auto-generated edit script applying numeric corrections to the Kujata_Profits workbook
# (per-row market price / profit recalculations produced by the scheduled data-refresh runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 334.10526
$ws.Range("I4").Value = 210.4375
$ws.Range("K4").Value = 210.4375
$ws.Range("M4").Value = -96.4375
$ws.Range("H40").Value = 2014.65
$ws.Range("J40").Value = 2081
$ws.Range("L40").Value = 2081
$ws.Range("N40").Value = -2431
$ws.Range("H61").Value = 158
$ws.Range("I61").Value = 120
$ws.Range("J61").Value = 500
$ws.Range("K61").Value = 360
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -188
$ws.Range("N61").Value = -1844
$ws.Range("H86").Value = 1876.3334
$ws.Range("I86").Value = 2255.0625
$ws.Range("K86").Value = 2255.0625
$ws.Range("M86").Value = -1132.0625
$ws.Range("H89").Value = 1876.3334
$ws.Range("I89").Value = 2255.0625
$ws.Range("K89").Value = 11275.3125
$ws.Range("M89").Value = -5659.3125
$ws.Range("H132").Value = 10108802
$ws.Range("I132").Value = 14500076
$ws.Range("J132").Value = 8871.799999999999
$ws.Range("K132").Value = 43500228
$ws.Range("L132").Value = 26615.4
$ws.Range("M132").Value = -43497698
$ws.Range("N132").Value = -31675.4
$ws.Range("H137").Value = 1422.2368
$ws.Range("I137").Value = 1075.762
$ws.Range("J137").Value = 1850.2354
$ws.Range("K137").Value = 3227.286
$ws.Range("L137").Value = 5550.706200000001
$ws.Range("M137").Value = -677.2860000000001
$ws.Range("N137").Value = -10650.7062
$ws.Range("H138").Value = 551953.4399999999
$ws.Range("I138").Value = 1292.0625
$ws.Range("J138").Value = 691804
$ws.Range("K138").Value = 3876.1875
$ws.Range("L138").Value = 2075412
$ws.Range("M138").Value = 1263.8125
$ws.Range("N138").Value = -2085692

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1209.1666
$ws.Range("I2").Value = 1126.375
$ws.Range("K2").Value = 1126.375
$ws.Range("M2").Value = -1013.375
$ws.Range("H61").Value = 142858620
$ws.Range("I61").Value = 200001070
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 200001070
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -200000858
$ws.Range("N61").Value = -2924
$ws.Range("H74").Value = 2388.5
$ws.Range("I74").Value = 1577
$ws.Range("K74").Value = 1577
$ws.Range("M74").Value = -703
$ws.Range("H77").Value = 2388.5
$ws.Range("I77").Value = 1577
$ws.Range("K77").Value = 7885
$ws.Range("M77").Value = -3517
$ws.Range("H116").Value = 1209.1666
$ws.Range("I116").Value = 1126.375
$ws.Range("K116").Value = 1126.375
$ws.Range("M116").Value = 1167.625
$ws.Range("H132").Value = 3327
$ws.Range("I132").Value = 3436
$ws.Range("K132").Value = 10308
$ws.Range("M132").Value = -7778
$ws.Range("H136").Value = 142858620
$ws.Range("I136").Value = 200001070
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 600003210
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -600000660
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1209.1666
$ws.Range("I3").Value = 1126.375
$ws.Range("K3").Value = 1126.375
$ws.Range("M3").Value = -1012.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1263.9348
$ws.Range("I31").Value = 1230.9111
$ws.Range("K31").Value = 1230.9111
$ws.Range("M31").Value = -935.9111
$ws.Range("H34").Value = 1263.9348
$ws.Range("I34").Value = 1230.9111
$ws.Range("K34").Value = 1230.9111
$ws.Range("M34").Value = -1028.9111
$ws.Range("H132").Value = 2060.5217
$ws.Range("I132").Value = 1885.5
$ws.Range("J132").Value = 2332.7778
$ws.Range("K132").Value = 5656.5
$ws.Range("L132").Value = 6998.3334
$ws.Range("M132").Value = -3126.5
$ws.Range("N132").Value = -12058.3334
$ws.Range("H134").Value = 15153161
$ws.Range("I134").Value = 1752.75
$ws.Range("J134").Value = 38463020
$ws.Range("K134").Value = 5258.25
$ws.Range("L134").Value = 115389060
$ws.Range("M134").Value = -2723.25
$ws.Range("N134").Value = -115394130

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 802.8570999999999
$ws.Range("I2").Value = 57.9
$ws.Range("J2").Value = 2665.25
$ws.Range("K2").Value = 347.4
$ws.Range("L2").Value = 15991.5
$ws.Range("M2").Value = -234.4
$ws.Range("N2").Value = -16217.5
$ws.Range("H56").Value = 6776.9165
$ws.Range("I56").Value = 6776.9165
$ws.Range("K56").Value = 6776.9165
$ws.Range("M56").Value = -6246.9165
$ws.Range("H122").Value = 1796.1177
$ws.Range("I122").Value = 850
$ws.Range("J122").Value = 1855.25
$ws.Range("K122").Value = 7650
$ws.Range("L122").Value = 16697.25
$ws.Range("M122").Value = -5200
$ws.Range("N122").Value = -21597.25
$ws.Range("H131").Value = 16950070
$ws.Range("J131").Value = 1074.3024
$ws.Range("L131").Value = 3222.9072
$ws.Range("N131").Value = -13302.9072

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 6600
$ws.Range("J109").Value = 6600
$ws.Range("L109").Value = 6600
$ws.Range("N109").Value = -8680
$ws.Range("H132").Value = 2746.7878
$ws.Range("I132").Value = 2481.8572
$ws.Range("K132").Value = 7445.571599999999
$ws.Range("M132").Value = -4915.571599999999
$ws.Range("H138").Value = 33503.168
$ws.Range("J138").Value = 33503.168
$ws.Range("L138").Value = 33503.168
$ws.Range("N138").Value = -43783.168

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3631.6924
$ws.Range("I40").Value = 2095.7778
$ws.Range("J40").Value = 7087.5
$ws.Range("K40").Value = 2095.7778
$ws.Range("L40").Value = 7087.5
$ws.Range("M40").Value = -1959.7778
$ws.Range("N40").Value = -7359.5
$ws.Range("H132").Value = 2722.8462
$ws.Range("I132").Value = 2492.6428
$ws.Range("J132").Value = 2991.4167
$ws.Range("K132").Value = 7477.928400000001
$ws.Range("L132").Value = 8974.250100000001
$ws.Range("M132").Value = -4947.928400000001
$ws.Range("N132").Value = -14034.2501
$ws.Range("H136").Value = 1497.25
$ws.Range("I136").Value = 996.3333
$ws.Range("K136").Value = 2988.9999
$ws.Range("M136").Value = -438.9998999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H120").Value = 32000
$ws.Range("J120").Value = 32000
$ws.Range("L120").Value = 32000
$ws.Range("N120").Value = -41676
$ws.Range("H132").Value = 2216.5454
$ws.Range("I132").Value = 1866.625
$ws.Range("K132").Value = 5599.875
$ws.Range("M132").Value = -3069.875

